# Fix duplicate "Anf:" label numbers in the requirements list.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Correct the duplicated ID values
$ws.Range("B17").Value = 4022
$ws.Range("B20").Value = 4040
$ws.Range("B21").Value = 4041
$ws.Range("B22").Value = 4042

# Update the visible window/selection to match where the edit was made
$ws.Activate()
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 13
